$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (price + 1h volume change) for the crypto snapshot.
# Each target cell stores its value as text (e.g. "307.49", "-4.70%"), so we
# force a Text number format before writing the new value, then restore the
# cell to the default "Normal" style so no spurious style diff is introduced.
function Set-TextValue($Sheet, $Addr, $Text) {
    $cell = $Sheet.Range($Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "307.49"
Set-TextValue $ws "E2" "-4.70%"
Set-TextValue $ws "D3" "40.02"
Set-TextValue $ws "E3" "-6.95%"
Set-TextValue $ws "D4" "5.132"
Set-TextValue $ws "E4" "-1.17%"
Set-TextValue $ws "D5" "0.07729"
Set-TextValue $ws "E5" "-5.87%"
Set-TextValue $ws "D6" "4.245"
Set-TextValue $ws "E6" "-1.70%"
Set-TextValue $ws "D7" "1.631"
Set-TextValue $ws "E7" "-10.56%"
Set-TextValue $ws "D8" "0.8808"
Set-TextValue $ws "E8" "-5.75%"
Set-TextValue $ws "D9" "0.1011"
Set-TextValue $ws "E9" "-9.01%"
Set-TextValue $ws "D10" "0.1751"
Set-TextValue $ws "E10" "-5.74%"
Set-TextValue $ws "D11" "0.09004"
Set-TextValue $ws "E11" "-4.44%"
Set-TextValue $ws "D12" "0.04391"
Set-TextValue $ws "E12" "-4.93%"
Set-TextValue $ws "D13" "0.1054"
Set-TextValue $ws "E13" "-0.14%"
Set-TextValue $ws "D14" "0.001260"
Set-TextValue $ws "E14" "-2.58%"
Set-TextValue $ws "D15" "0.005861"
Set-TextValue $ws "E15" "0.94%"
Set-TextValue $ws "D16" "3.354"
Set-TextValue $ws "D17" "2.423"
Set-TextValue $ws "E17" "-4.75%"
Set-TextValue $ws "D18" "0.3324"
Set-TextValue $ws "E18" "-0.47%"
Set-TextValue $ws "D19" "7.047"
Set-TextValue $ws "E19" "-4.51%"
Set-TextValue $ws "E20" "-3.46%"
Set-TextValue $ws "D21" "0.2995"
Set-TextValue $ws "E21" "14.37%"
Set-TextValue $ws "D22" "0.04166"
Set-TextValue $ws "E22" "0.10%"
Set-TextValue $ws "E23" "-3.48%"
Set-TextValue $ws "D24" "0.004119"
Set-TextValue $ws "E24" "-4.99%"
Set-TextValue $ws "D25" "0.0001300"
Set-TextValue $ws "E25" "8.46%"
Set-TextValue $ws "E26" "0.23%"
Set-TextValue $ws "D38" "0.02361"
Set-TextValue $ws "E38" "-14.65%"
Set-TextValue $ws "D39" "0.05165"
Set-TextValue $ws "E39" "-7.66%"
Set-TextValue $ws "D40" "0.007919"
Set-TextValue $ws "E40" "-0.91%"
Set-TextValue $ws "E41" "-5.03%"
Set-TextValue $ws "D42" "0.006375"
Set-TextValue $ws "E42" "-2.32%"
Set-TextValue $ws "D43" "0.001963"
Set-TextValue $ws "E43" "-5.72%"
Set-TextValue $ws "D44" "0.008489"
Set-TextValue $ws "E44" "13.23%"
Set-TextValue $ws "D45" "0.3316"
Set-TextValue $ws "E45" "-6.43%"
Set-TextValue $ws "D46" "0.00006522"
Set-TextValue $ws "E46" "-6.62%"
Set-TextValue $ws "E47" "0.12%"
Set-TextValue $ws "E48" "98.76%"
Set-TextValue $ws "D49" "0.002177"
Set-TextValue $ws "E49" "-37.46%"
Set-TextValue $ws "E50" "0.12%"
Set-TextValue $ws "E51" "0.12%"
